# Reorder the "Recorded By" (column G) entries so that known user-email
# addresses are listed before the "System"/"system" entries, leaving the
# relative order of the remaining entries untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The specific addresses that get promoted to the front of the list.
$priorityEmails = @("dnasr281@gmail.com", "backup@backdoor.com")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }
    if (-not $current.Contains(",")) {
        continue
    }

    $parts = $current.Split(",")
    $priorityParts = @()
    $otherParts = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($priorityEmails -contains $trimmed) {
            $priorityParts += $trimmed
        } else {
            $otherParts += $trimmed
        }
    }

    if ($priorityParts.Count -eq 0) {
        continue
    }

    $newParts = $priorityParts + $otherParts
    $newValue = $newParts -join ", "

    if ($newValue -ne $current) {
        $cell.Value2 = $newValue
    }
}
